# Commit: "added common method to set status; added test for adding new centers"
#
# - rename Sheet2 -> add_new_centers, and populate it with a new
#   "add_new_centers" test-data table (mirrors add_new_locations/
#   add_institute_details)
# - add_new_locations (sheet1): the location "code" values used by the
#   tests were regenerated (loc001/loc002/loc003/branch1/bra2/brrr3/b333
#   -> loc999/loczbc/locdfdfdf/loc093/loc34343/loc---/locxxxx)
# - active tab moves from Sheet2 to add_new_locations

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the old throwaway "Sheet2" to "add_new_centers"
# ---------------------------------------------------------------------
$wsCenters = $wb.Worksheets.Item(3)
$wsCenters.Name = "add_new_centers"

# ---------------------------------------------------------------------
# 2. add_new_locations: refresh the generated "code" column values
# ---------------------------------------------------------------------
$wsLocations = $wb.Worksheets.Item(1)

$wsLocations.Cells.Item(2, 1).Value = "loc999"
$wsLocations.Cells.Item(3, 1).Value = "loczbc"
$wsLocations.Cells.Item(4, 1).Value = "locdfdfdf"
$wsLocations.Cells.Item(5, 1).Value = "loc093"
$wsLocations.Cells.Item(6, 1).Value = "loc34343"
$wsLocations.Cells.Item(7, 1).Value = "loc---"
$wsLocations.Cells.Item(8, 1).Value = "locxxxx"

# ---------------------------------------------------------------------
# 3. add_new_centers: fill in the new test-data table
# ---------------------------------------------------------------------
$wsCenters.Cells.Item(1, 1).Value = "code"
$wsCenters.Cells.Item(1, 2).Value = "name"
$wsCenters.Cells.Item(1, 3).Value = "location"
$wsCenters.Cells.Item(1, 4).Value = "status"
$wsCenters.Cells.Item(1, 5).Value = "runmode"

$wsCenters.Cells.Item(2, 1).Value = "CENx"
$wsCenters.Cells.Item(2, 2).Value = "exam center"
$wsCenters.Cells.Item(2, 3).Value = "colombo"
$wsCenters.Cells.Item(2, 4).Value = "Active"
$wsCenters.Cells.Item(2, 5).Value = "Y"

$wsCenters.Cells.Item(3, 1).Value = "CENy"
$wsCenters.Cells.Item(3, 2).Value = "exam center"
$wsCenters.Cells.Item(3, 3).Value = "metro"
$wsCenters.Cells.Item(3, 4).Value = "Inactive"
$wsCenters.Cells.Item(3, 5).Value = "Y"

$wsCenters.Cells.Item(4, 1).Value = "CENi"
$wsCenters.Cells.Item(4, 2).Value = "library"
$wsCenters.Cells.Item(4, 3).Value = "malabe"
$wsCenters.Cells.Item(4, 4).Value = "Active"
$wsCenters.Cells.Item(4, 5).Value = "Y"

$wsCenters.Cells.Item(5, 1).Value = "CEN1"
$wsCenters.Cells.Item(5, 2).Value = "canteen"
$wsCenters.Cells.Item(5, 3).Value = "jaela"
$wsCenters.Cells.Item(5, 4).Value = "Inactive"
$wsCenters.Cells.Item(5, 5).Value = "Y"

$wsCenters.Cells.Item(6, 1).Value = "CEN-"
$wsCenters.Cells.Item(6, 2).Value = "recreation area"
$wsCenters.Cells.Item(6, 3).Value = "kandy"
$wsCenters.Cells.Item(6, 4).Value = "Inactive"
$wsCenters.Cells.Item(6, 5).Value = "Y"

$wsCenters.Cells.Item(7, 1).Value = "CEN00"
$wsCenters.Cells.Item(7, 2).Value = "wash room"
$wsCenters.Cells.Item(7, 3).Value = "malabe"
$wsCenters.Cells.Item(7, 4).Value = "Active"
$wsCenters.Cells.Item(7, 5).Value = "Y"

# keep the pre-existing selection rectangle (A8) on the newly-built sheet
$wsCenters.Range("A8").Select()

# ---------------------------------------------------------------------
# 4. Move the active tab / selection to add_new_locations
# ---------------------------------------------------------------------
$wsLocations.Range("A9").Select()
